# Add the new "Rusing on paper clips..." daily log entry (row 7) on the
# "Daily" sheet, then leave the Daily sheet active/selected at E8 (matching
# where Excel would land right after typing the last cell of the new row).
# This also flips "Weekly"'s tabSelected off, since only one sheet can be
# the active tab in a workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily")

$ws.Cells.Item(7, 1).Value = 20221004
$ws.Cells.Item(7, 2).Value = "completed"
$ws.Cells.Item(7, 3).Value = "completed"
$ws.Cells.Item(7, 4).Value = "AH"
$ws.Cells.Item(7, 5).Value = "Rusing on paper clips, will need to replace soon"

$ws.Activate()
$ws.Range("E8").Select()
